# Apply numeric updates to the "F" (人气值/热度 or similar) column across sheets
# as described by the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 537
$wsExhibition.Range("F4").Value = 47
$wsExhibition.Range("F9").Value = 380
$wsExhibition.Range("F10").Value = 3412

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 95

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 95
$wsAll.Range("F4").Value = 537
$wsAll.Range("F5").Value = 47
$wsAll.Range("F10").Value = 380
$wsAll.Range("F11").Value = 3412
